$d = $word.ActiveDocument

# Each pair is (old text, new text). Every "old" string occurs exactly once
# in the original document, so ReplaceAll (wdReplaceAll = 2) on the whole
# document Content range safely retargets just that single occurrence, even
# though one "new" value later coincides with an "old" value from an earlier
# step (handled safely because we search in the original uniqueness order).
$pairs = @(
    @("2023-12-03 Sunday", "2023-12-04 Monday"),
    @("53×45=2385", "72×29=2088"),
    @("57×93=5301", "93×58=5394"),
    @("48×42=2016", "79×13=1027"),
    @("36×83=2988", "76×72=5472"),
    @("37×30=1110", "80×50=4000"),
    @("91×75=6825", "43×69=2967"),
    @("72×30=2160", "58×72=4176"),
    @("63×96=6048", "47×57=2679"),
    @("98×81=7938", "36×83=2988"),
    @("18×56=1008", "71×43=3053"),
    @("80×91=7280", "45×11=495"),
    @("48×76=3648", "69×72=4968"),
    @("71×58=4118", "96×19=1824"),
    @("69×99=6831", "91×17=1547"),
    @("96×64=6144", "88×26=2288"),
    @("58×91=5278", "34×58=1972"),
    @("15×16=240", "97×13=1261"),
    @("17×76=1292", "99×89=8811"),
    @("67×44=2948", "98×93=9114"),
    @("46×71=3266", "99×77=7623"),
    @("44×74=3256", "48×59=2832"),
    @("21×87=1827", "19×84=1596"),
    @("29×50=1450", "42×97=4074"),
    @("17×36=612", "55×91=5005"),
    @("63×75=4725", "59×29=1711")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: replacement not found for '$old'"
    }
}

Write-Host "Replacements complete"
